$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: round Q/R to nearest integer, clear Starttid/Sluttid (Z/AB)
$ws.Range("Q2").Value = 547961
$ws.Range("R2").Value = 6960421
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# Row 3: swap Id with row 4, update Q/R, clear Starttid/Sluttid
$ws.Range("A3").Value = 112043819
$ws.Range("Q3").Value = 547979
$ws.Range("R3").Value = 6960195
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# Row 4: swap Id with row 3, update Q/R, clear Starttid/Sluttid
$ws.Range("A4").Value = 112043839
$ws.Range("Q4").Value = 547969
$ws.Range("R4").Value = 6960405
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
